$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("E2").Value = [double]"24.76000000000043"
$ws.Range("H2").Value = [double]"1.329608412724738e-16"
$ws.Range("K2").Value = [double]"59.6823584242938"
$ws.Range("L2").Value = "[53.81119771729911, 65.55351913128848]"
$ws.Range("O2").Value = [double]"1.578658170272348"
$ws.Range("P2").Value = "[1.46544762419704, 1.6918687163476562]"
$ws.Range("S2").Value = [double]"54.34120132165868"
$ws.Range("T2").Value = "[50.30761438585791, 58.37478825745945]"
$ws.Range("W2").Value = [double]"18.53901901901934"
$ws.Range("X2").Value = [double]"18.09289289289321"
$ws.Range("Y2").Value = [double]"18.98514514514548"

# Row 3
$ws.Range("E3").Value = [double]"25.78000000000059"
$ws.Range("H3").Value = [double]"1.329608412724738e-16"
$ws.Range("K3").Value = [double]"59.03825967883289"
$ws.Range("L3").Value = "[52.77603306738779, 65.300486290278]"
$ws.Range("O3").Value = [double]"2.496921488438735"
$ws.Range("P3").Value = "[2.3837109423634257, 2.6101320345140437]"
$ws.Range("S3").Value = [double]"54.74759334609768"
$ws.Range("T3").Value = "[51.34661140654389, 58.14857528565147]"
$ws.Range("W3").Value = [double]"15.53509509509545"
$ws.Range("X3").Value = [double]"15.07059059059093"
$ws.Range("Y3").Value = [double]"15.99959959959997"
